# Merge start: 1st Row1, other Row2
# Insert a new header row above the existing data (all data rows shift down
# by one: old row 1 -> row 2, ... old row 6 -> row 7) and populate the new
# header row with "Time" / "kind" / "num".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything down by inserting a fresh row at the top.
$ws.Rows.Item(1).Insert()

# Touch alignment (re-asserting the existing default) so the new header
# cells pick up their own style entry, then write the header labels.
$ws.Range("A1:C1").VerticalAlignment = -4107
$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "kind"
$ws.Range("C1").Value = "num"

# Leave the new header row selected.
$ws.Range("A1:C1").Select() | Out-Null
